# Generate Report for Handoff
# - Flip status from "In Translation" to "Ready for handoff" on all three
#   sheets (Overview zh-cn/de-de columns, and the Status column on each
#   per-locale sheet).
# - Bump the associated "Latest HO Xliff Generate Date" / "Latest Handoff
#   Datetime" timestamps forward by ~70s/~80s to the new handoff moment.
# - Widen the affected status/date columns to fit the new, longer text.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# New column width that fits "Ready for handoff"/the widened date columns.
# (ColumnWidth is quantized to the sheet's character grid on write, so we
# feed it the value that lands closest to the target width after rounding.)
$newColWidth = 16.3

# --- Overview sheet ---------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = $newStatus            # zh-cn status
$overview.Range("F2").Value = $newStatus            # de-de status
$overview.Range("G2").Value = "2016-10-19 17:33:58" # Latest HO Xliff Generate Date

$overview.Columns.Item(5).ColumnWidth = $newColWidth # zh-cn column
$overview.Columns.Item(6).ColumnWidth = $newColWidth # de-de column

# --- zh-cn sheet --------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $newStatus                # Status
$zhcn.Range("H2").Value = "2016-10-19 17:33:47"     # Latest Handoff Datetime

$zhcn.Columns.Item(3).ColumnWidth = $newColWidth     # Status column

# --- de-de sheet --------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $newStatus                # Status
$dede.Range("H2").Value = "2016-10-19 17:33:58"     # Latest Handoff Datetime

$dede.Columns.Item(3).ColumnWidth = $newColWidth     # Status column
